# The deck ends with a run of "Reply to:" template slides (#2..#5 at
# slides 7-10). This edit trims the last two of them (#4 and #5, i.e.
# slides 9 and 10) along with their associated notes pages, which
# PowerPoint removes automatically when the owning slide is deleted.

$p = $ppt.ActivePresentation

# Delete from the end so indices of the remaining slides don't shift
# out from under us.
$p.Slides.Item($p.Slides.Count).Delete()
$p.Slides.Item($p.Slides.Count).Delete()
